$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: set species to Fisher and male count to 1
$ws.Range("F2").Value = "Fisher"
$ws.Range("H2").Value = 1

# Remove row 3 entirely (its data has been merged into row 2)
$ws.Rows(3).Delete()
